$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.923999999999999
$ws.Range("B3").Value = 7.7509999999999994
$ws.Range("B4").Value = 0.81800000000000006
$ws.Range("B6").Value = 7.1899999999999995
$ws.Range("B7").Value = 25.033999999999999
$ws.Range("B8").Value = 0.6
$ws.Range("B9").Value = 0
